$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text format so that Excel does not
# auto-convert numeric-looking strings (e.g. "593.21") into numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "66.855.45"
$ws.Range("E2").Value = "  -4.45%  "
$ws.Range("D3").Value = "3.221.61"
$ws.Range("E3").Value = "  -7.91%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "593.21"
$ws.Range("E5").Value = "  -1.92%  "
$ws.Range("D6").Value = "151.88"
$ws.Range("E6").Value = "  -11.80%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "3.212.25"
$ws.Range("E8").Value = "  -7.93%  "
$ws.Range("D9").Value = "0.546"
$ws.Range("E9").Value = "  -10.07%  "
$ws.Range("D10").Value = "0.174"
$ws.Range("E10").Value = "  -10.02%  "
$ws.Range("D11").Value = "6.66"
$ws.Range("E11").Value = "  -7.66%  "
$ws.Range("D12").Value = "0.502"
$ws.Range("E12").Value = "  -14.08%  "
$ws.Range("D13").Value = "39.23"
$ws.Range("E13").Value = "  -14.59%  "
$ws.Range("D14").Value = "0.0000246"
$ws.Range("E14").Value = "  -10.28%  "
$ws.Range("D15").Value = "3.740.34"
$ws.Range("E15").Value = "  -8.05%  "
$ws.Range("D16").Value = "66.889.50"
$ws.Range("E16").Value = "  -4.49%  "
$ws.Range("D17").Value = "3.215.90"
$ws.Range("E17").Value = "  -8.00%  "
$ws.Range("D18").Value = "0.115"
$ws.Range("E18").Value = "  -4.80%  "
$ws.Range("D19").Value = "7.22"
$ws.Range("E19").Value = "  -13.45%  "
$ws.Range("D20").Value = "532.34"
$ws.Range("E20").Value = "  -13.09%  "
$ws.Range("D21").Value = "15.09"
$ws.Range("E21").Value = "  -13.60%  "
$ws.Range("D22").Value = "0.764"
$ws.Range("E22").Value = "  -12.77%  "
$ws.Range("D23").Value = "7.98"
$ws.Range("E23").Value = "  -12.46%  "
$ws.Range("D24").Value = "13.86"
$ws.Range("E24").Value = "  -10.42%  "
$ws.Range("D25").Value = "86.21"
$ws.Range("E25").Value = "  -12.42%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "3.19"
$ws.Range("E27").Value = "  -14.15%  "
$ws.Range("D28").Value = "2.21"
$ws.Range("E28").Value = "  -13.31%  "
$ws.Range("D29").Value = "8.16"
$ws.Range("E29").Value = "  -8.90%  "
$ws.Range("D30").Value = "29.47"
$ws.Range("E30").Value = "  -12.51%  "
$ws.Range("D31").Value = "2.64"
$ws.Range("E31").Value = "  -11.24%  "
$ws.Range("D32").Value = "1.14"
$ws.Range("E32").Value = "  -10.59%  "
$ws.Range("D33").Value = "543.03"
$ws.Range("E33").Value = "  -13.43%  "
$ws.Range("D34").Value = "6.58"
$ws.Range("E34").Value = "  -17.97%  "
$ws.Range("D35").Value = "5.76"
$ws.Range("E35").Value = "  -15.07%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").Value = "53.13"
$ws.Range("E37").Value = "  -6.19%  "
$ws.Range("D38").Value = "0.0874"
$ws.Range("E38").Value = "  -12.07%  "
$ws.Range("D39").Value = "9.37"
$ws.Range("E39").Value = "  -12.48%  "
$ws.Range("D40").Value = "0.0423"
$ws.Range("E40").Value = "  -11.74%  "
$ws.Range("E41").Value = "  -12.33%  "
$ws.Range("D42").Value = "2.939.74"
$ws.Range("E42").Value = "  -12.20%  "
$ws.Range("D43").Value = "2.64"
$ws.Range("E43").Value = "  -23.15%  "
$ws.Range("D44").Value = "0.266"
$ws.Range("E44").Value = "  -13.75%  "
$ws.Range("D45").Value = "0.0₃0588"
$ws.Range("E45").Value = "  -19.30%  "
$ws.Range("D46").Value = "2.42"
$ws.Range("E46").Value = "  -16.16%  "
$ws.Range("D47").Value = "26.54"
$ws.Range("E47").Value = "  -16.33%  "
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("E49").Value = "  -16.42%  "
$ws.Range("D50").Value = "0.115"
$ws.Range("E50").Value = "  -11.56%  "
$ws.Range("D51").Value = "122.74"
$ws.Range("E51").Value = "  -7.75%  "

# Restore the original (default/Normal) cell style so the saved file
# matches the original formatting of the worksheet.
$dataRange.Style = "Normal"
